$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled dSF (column F) data for 2021 rows - mean calculation update
$ws.Range("F2").Value = 3
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -3
$ws.Range("F15").Value = -4
$ws.Range("F16").Value = -5
$ws.Range("F18").Value = -5
$ws.Range("F19").Value = 7
